$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update indicator information block (row 4) and data-reporter block (rows 6-10)
$ws.Range("B4").Value = "6.2.1 Proportion of population using (a) safely managed sanitation services and (b) a hand-washing facility with soap and water"
$ws.Range("B6").Value = "National Statistical Committee of the Kyrgyz Republic (Department of Household Statistics)"
$ws.Range("B7").Value = "Kalymbetova Yryskan"
$ws.Range("B8").Value = "yryskan.kalymbetova@gmail.com "
$ws.Range("B9").Value = "(0312) 32 46 55"
$ws.Range("B10").Value = "www.stat.gov.kg"

# Move the sheet selection from B2 to B4
$ws.Range("B4").Select()
